$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("credentials")

# ---------------------------------------------------------------------------
# Sheet1: populate the previously-blank grid with the credentials data
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value2 = "userName"
$ws1.Range("B1").Value2 = "password"
$ws1.Range("C1").Value2 = "employeeToAdd"
$ws1.Range("D1").Value2 = "newUserName"

$ws1.Range("A2").Value2 = "Admin"
$ws1.Range("B2").Value2 = "OYs6MbnC2@"
$ws1.Range("C2").Value2 = "Aaron"
$ws1.Range("D2").Value2 = "Aaron Update"

$ws1.Range("A3").Value2 = "Vish"
$ws1.Range("C3").Value2 = "Vishwanath"
$ws1.Range("D3").Value2 = "Vishwanath D B"

# B3 gets a mailto hyperlink (matches the "credentials" sheet's existing one);
# add the hyperlink first then restore the plain display text of the cell so
# the visible value stays "Password@123" while the link's legacy display
# text is the mailto address.
$ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:Password@123", [Type]::Missing, [Type]::Missing, "mailto:Password@123")
$ws1.Range("B3").Value2 = "Password@123"

# ---------------------------------------------------------------------------
# credentials: insert a new second row (Tabby / Listener+Reporting creds)
# above the existing "Admin / OYs6MbnC2@ / Aaron / Aaron Update" row.
# ---------------------------------------------------------------------------
$ws2.Rows("2:2").Insert()

$ws2.Range("A2").Value2 = "Admin"
$ws2.Range("A2").Style = "Normal"

$ws2.Range("C2").Value2 = "Tabby"
$ws2.Range("C2").Style = "Normal"

$ws2.Range("D2").Value2 = "TabbyUpdate"
$ws2.Range("D2").Style = "Normal"

$ws2.Range("B2").Value2 = "iVLq@23JjQ"
$ws2.Range("B2").Style = "Normal"
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:iVLq@23JjQ")
$ws2.Range("B2").Value2 = "iVLq@23JjQ"

# ---------------------------------------------------------------------------
# Selections: matches the saved cursor position on each sheet.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A2:D3").Select()

$ws2.Activate()
$ws2.Range("A3:D4").Select()

Write-Output "done"
